$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 is North Carolina - fill in the results from the July 4, 2020 run.

# B4: Date Published - same numeric serial + style as the other date cells (e.g. B2/B3)
$ws.Range("B4").Value = 44016
$ws.Range("B4").Style = $ws.Range("B2").Style
$ws.Range("B4").NumberFormat = $ws.Range("B2").NumberFormat

$ws.Range("C4").Value = 71654
$ws.Range("D4").Value = 1395
$ws.Range("E4").Value = 11390
$ws.Range("F4").Value = 446
$ws.Range("G4").Value = 22.98
$ws.Range("H4").Value = 33.16

# I4 (Pct Includes Unknown Race) stays FALSE; J4 (Pct Includes Hispanic Black) flips to TRUE
$ws.Range("J4").Value = $true

$ws.Range("K4").Value = 49566
$ws.Range("L4").Value = 1345

# O4: Status code goes from the timeout error message to Success!
$ws.Range("O4").Value = "Success!"
